$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166277408599854
$ws.Range("B1").Value = 2.432313203811646
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.370549917221069
$ws.Range("E1").Value = 1.235013604164124
